# Regional Availability Factor.xlsx - "updated 4.0 files and mdl"
#
# 1. About sheet: bump the "last updated" date (C1) from 2024-03-15 to 2024-03-28
# 2. RAF-capacity sheet: raise the capacity-credit RAF for the two hydrogen
#    technologies (rows 24 & 25, column B) from 0.3 to 1
# 3. Move the active tab / selection from RAF-generation to RAF-capacity,
#    leaving the cursor on B25 and zoomed to 80%.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")

# --- About: refresh the date stamp in C1 ---
$wsAbout.Range("C1").Value = 45379

# --- RAF-capacity: hydrogen combustion turbine / hydrogen combined cycle RAF ---
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# --- View state: make RAF-capacity the active/selected sheet ---
$wsCapacity.Activate()
[void]$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.Zoom = 80
